$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.107667207717896
$ws.Range("B1").Value = 2.419713735580444
$ws.Range("C1").Value = 2.079732656478882
$ws.Range("D1").Value = 2.205829858779907
$ws.Range("E1").Value = 2.630209445953369
